# Apply the "create insertion query from sheet must give value types" edit.
#
# Summary of the change:
#  - Sheet "Sheet2" (first tab) is rewritten from a tiny 2x2 "Hello/Dude" demo
#    table into a one-row "dobpaymentrequest" table: a header row naming the
#    DB columns of a payment-request insert, and a data row holding the
#    values for that insert - numeric columns stay numeric, text/date/json
#    columns are single-quoted SQL string literals.
#  - Sheet "Sheet1" (second tab) keeps its existing fruit-table content; only
#    its shared-string indices shift because of the strings inserted ahead of
#    them in sharedStrings.xml, plus a zoom-level bump.
#  - Workbook-level tab ratio / window zoom tweaks.

$wb = $excel.ActiveWorkbook

# Narrower sheet-tab area / wider horizontal scrollbar in the window chrome.
$excel.ActiveWindow.TabRatio = 0.5

$ws1 = $wb.Worksheets.Item("Sheet2")
$ws2 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Sheet2: rebuild as the dobpaymentrequest insert table
# ---------------------------------------------------------------------
[void]$ws1.Activate()

# Clear the old 2x2 demo content first.
[void]$ws1.Cells.Clear()

$ws1.Range("A1").Value = "dobpaymentrequest"

$ws1.Range("A2").Value = "id"
$ws1.Range("B2").Value = "code"
$ws1.Range("C2").Value = "creationdate"
$ws1.Range("D2").Value = "modifieddate"
$ws1.Range("E2").Value = "creationinfo"
$ws1.Range("F2").Value = "modificationinfo"
$ws1.Range("G2").Value = "currentstates"
$ws1.Range("H2").Value = "paymentType"
$ws1.Range("I2").Value = "purchaseUnitId"
$ws1.Range("J2").Value = "paymentForm"
$ws1.Range("K2").Value = "amount"

$ws1.Range("A3").Value = 1001
$ws1.Range("B3").Value = "''2019000001'"
$ws1.Range("C3").Value = "''2018-08-05 09:02:00'"
$ws1.Range("D3").Value = "''2018-08-05 09:02:00'"
$ws1.Range("E3").Value = "''Amira.Atya'"
$ws1.Range("F3").Value = "''Amira.Atya'"
$ws1.Range("G3").Value = "''[""Draft""]'"
$ws1.Range("H3").Value = "''GENERAL'"
$ws1.Range("I3").Value = 20
$ws1.Range("J3").Value = "''CASH" + [char]0x2019
$ws1.Range("K3").Value = 1000

# Column widths (characters).
$ws1.Columns.Item(1).ColumnWidth = 17.5
$ws1.Columns.Item(2).ColumnWidth = 12.17
$ws1.Columns.Item(3).ColumnWidth = 18.33
$ws1.Columns.Item(4).ColumnWidth = 22.42
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 15.33
$ws1.Columns.Item(8).ColumnWidth = 12.17
$ws1.Columns.Item(9).ColumnWidth = 13.76
$ws1.Columns.Item(10).ColumnWidth = 12.83
$ws1.Columns.Item(11).ColumnWidth = 9.67

$excel.ActiveWindow.Zoom = 120
[void]$ws1.Range("H8").Select()

# ---------------------------------------------------------------------
# Sheet1: unchanged content, just the zoom bump (shared-string indices
# shift automatically because of the new strings added ahead of them).
# ---------------------------------------------------------------------
[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 120
[void]$ws2.Range("G11").Select()

# Restore the originally active sheet/tab.
[void]$ws1.Activate()
